# Changed date and time fields for JGI app and verified database persistence
#
# The "survey" sheet's `type` column (C) listed "date" for the FOL_date
# question and "time" for the FOL_time_begin / FOL_time_end questions.
# Those field types are changed to plain "text" fields.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")

# FOL_date (row 2): type date -> text
$ws.Range("C2").Value = "text"

# FOL_time_begin (row 5): type time -> text
$ws.Range("C5").Value = "text"

# FOL_time_end (row 6): type time -> text
$ws.Range("C6").Value = "text"

# Update the active selection on the survey sheet to reflect where the
# author ended up after making the edit.
$ws.Range("C7").Select() | Out-Null
